$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '30.365.42'
$ws.Range('D2').Style = 'Normal'

$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +2.62%  '
$ws.Range('E2').Style = 'Normal'

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.108.47'
$ws.Range('D3').Style = 'Normal'

$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +0.77%  '
$ws.Range('E3').Style = 'Normal'

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.005'
$ws.Range('D4').Style = 'Normal'

$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.42%  '
$ws.Range('E4').Style = 'Normal'

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '344.52'
$ws.Range('D5').Style = 'Normal'

$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +0.82%  '
$ws.Range('E5').Style = 'Normal'

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.004'
$ws.Range('D6').Style = 'Normal'

$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -0.32%  '
$ws.Range('E6').Style = 'Normal'

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5237'
$ws.Range('D7').Style = 'Normal'

$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +2.23%  '
$ws.Range('E7').Style = 'Normal'

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.4444'
$ws.Range('D8').Style = 'Normal'

$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +1.18%  '
$ws.Range('E8').Style = 'Normal'

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '54.64'
$ws.Range('D9').Style = 'Normal'

$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +2.93%  '
$ws.Range('E9').Style = 'Normal'

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.09463'
$ws.Range('D10').Style = 'Normal'

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.174'
$ws.Range('D11').Style = 'Normal'

$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +0.62%  '
$ws.Range('E11').Style = 'Normal'

$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +1.83%  '
$ws.Range('E12').Style = 'Normal'

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '8.729'
$ws.Range('D13').Style = 'Normal'

$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +6.59%  '
$ws.Range('E13').Style = 'Normal'

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.953'
$ws.Range('D14').Style = 'Normal'

$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +2.99%  '
$ws.Range('E14').Style = 'Normal'

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '2.079.86'
$ws.Range('D15').Style = 'Normal'

$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -0.98%  '
$ws.Range('E15').Style = 'Normal'

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '101.68'
$ws.Range('D16').Style = 'Normal'

$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +2.02%  '
$ws.Range('E16').Style = 'Normal'

$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +1.83%  '
$ws.Range('E17').Style = 'Normal'

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '1.006'
$ws.Range('D18').Style = 'Normal'

$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -0.32%  '
$ws.Range('E18').Style = 'Normal'

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '21.28'
$ws.Range('D19').Style = 'Normal'

$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +1.42%  '
$ws.Range('E19').Style = 'Normal'

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.06716'
$ws.Range('D20').Style = 'Normal'

$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +1.21%  '
$ws.Range('E20').Style = 'Normal'

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.351'
$ws.Range('D21').Style = 'Normal'

$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +3.02%  '
$ws.Range('E21').Style = 'Normal'

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.004'
$ws.Range('D22').Style = 'Normal'

$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -0.41%  '
$ws.Range('E22').Style = 'Normal'

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '30.389.87'
$ws.Range('D23').Style = 'Normal'

$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +2.53%  '
$ws.Range('E23').Style = 'Normal'

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '12.66'
$ws.Range('D24').Style = 'Normal'

$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +0.67%  '
$ws.Range('E24').Style = 'Normal'

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.309'
$ws.Range('D25').Style = 'Normal'

$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +0.12%  '
$ws.Range('E25').Style = 'Normal'

$ws.Range('B26').Value = 'EthereumClassic'

$ws.Range('C26').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '22.03'
$ws.Range('D26').Style = 'Normal'

$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +1.29%  '
$ws.Range('E26').Style = 'Normal'

$ws.Range('B27').Value = 'Monero'

$ws.Range('C27').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '163.53'
$ws.Range('D27').Style = 'Normal'

$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +0.68%  '
$ws.Range('E27').Style = 'Normal'

$ws.Range('B28').Value = 'LidoDAOToken'

$ws.Range('C28').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.543'
$ws.Range('D28').Style = 'Normal'

$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +0.92%  '
$ws.Range('E28').Style = 'Normal'

$ws.Range('B29').Value = 'BitcoinCash'

$ws.Range('C29').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '133.85'
$ws.Range('D29').Style = 'Normal'

$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +1.22%  '
$ws.Range('E29').Style = 'Normal'

$ws.Range('B30').Value = 'ImmutableX'

$ws.Range('C30').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.154'
$ws.Range('D30').Style = 'Normal'

$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +2.50%  '
$ws.Range('E30').Style = 'Normal'

$ws.Range('B31').Value = 'ARBITRUM'

$ws.Range('C31').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.748'
$ws.Range('D31').Style = 'Normal'

$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +7.25%  '
$ws.Range('E31').Style = 'Normal'

$ws.Range('B32').Value = 'Stellar'

$ws.Range('C32').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.1054'
$ws.Range('D32').Style = 'Normal'

$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +1.09%  '
$ws.Range('E32').Style = 'Normal'

$ws.Range('B33').Value = 'InternetComputer(DFINITY)'

$ws.Range('C33').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '6.865'
$ws.Range('D33').Style = 'Normal'

$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +13.69%  '
$ws.Range('E33').Style = 'Normal'

$ws.Range('B34').Value = 'Filecoin'

$ws.Range('C34').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '6.275'
$ws.Range('D34').Style = 'Normal'

$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +2.40%  '
$ws.Range('E34').Style = 'Normal'

$ws.Range('B35').Value = 'HuobiToken'

$ws.Range('C35').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.923'
$ws.Range('D35').Style = 'Normal'

$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -0.82%  '
$ws.Range('E35').Style = 'Normal'

$ws.Range('B36').Value = 'FraxShare'

$ws.Range('C36').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '10.45'
$ws.Range('D36').Style = 'Normal'

$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +2.04%  '
$ws.Range('E36').Style = 'Normal'

$ws.Range('B37').Value = 'VeChain'

$ws.Range('C37').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.02630'
$ws.Range('D37').Style = 'Normal'

$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +2.61%  '
$ws.Range('E37').Style = 'Normal'

$ws.Range('B38').Value = 'Hedera'

$ws.Range('C38').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.06811'
$ws.Range('D38').Style = 'Normal'

$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +2.37%  '
$ws.Range('E38').Style = 'Normal'

$ws.Range('B39').Value = 'TheSandbox'

$ws.Range('C39').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.7068'
$ws.Range('D39').Style = 'Normal'

$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +3.56%  '
$ws.Range('E39').Style = 'Normal'

$ws.Range('B40').Value = 'Aptos'

$ws.Range('C40').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '12.58'
$ws.Range('D40').Style = 'Normal'

$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +1.94%  '
$ws.Range('E40').Style = 'Normal'

$ws.Range('B41').Value = 'TrustWalletToken'

$ws.Range('C41').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.347'
$ws.Range('D41').Style = 'Normal'

$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +5.45%  '
$ws.Range('E41').Style = 'Normal'

$ws.Range('B42').Value = 'Algorand'

$ws.Range('C42').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.2230'
$ws.Range('D42').Style = 'Normal'

$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +0.28%  '
$ws.Range('E42').Style = 'Normal'

$ws.Range('B43').Value = 'Decentraland'

$ws.Range('C43').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.6854'
$ws.Range('D43').Style = 'Normal'

$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +3.40%  '
$ws.Range('E43').Style = 'Normal'

$ws.Range('B44').Value = 'EnergySwap'

$ws.Range('C44').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '14.46'
$ws.Range('D44').Style = 'Normal'

$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +2.72%  '
$ws.Range('E44').Style = 'Normal'

$ws.Range('B45').Value = 'NEARProtocol'

$ws.Range('C45').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.368'
$ws.Range('D45').Style = 'Normal'

$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +3.59%  '
$ws.Range('E45').Style = 'Normal'

$ws.Range('B46').Value = 'Frax'

$ws.Range('C46').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.004'
$ws.Range('D46').Style = 'Normal'

$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -0.28%  '
$ws.Range('E46').Style = 'Normal'

$ws.Range('B47').Value = 'WEMIXTOKEN'

$ws.Range('C47').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.385'
$ws.Range('D47').Style = 'Normal'

$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +18.92%  '
$ws.Range('E47').Style = 'Normal'

$ws.Range('B48').Value = 'PancakeSwap'

$ws.Range('C48').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.654'
$ws.Range('D48').Style = 'Normal'

$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +1.40%  '
$ws.Range('E48').Style = 'Normal'

$ws.Range('B49').Value = 'BabyDogeCoin'

$ws.Range('C49').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.00000000348'
$ws.Range('D49').Style = 'Normal'

$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +3.71%  '
$ws.Range('E49').Style = 'Normal'

$ws.Range('B50').Value = 'ThetaToken'

$ws.Range('C50').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.207'
$ws.Range('D50').Style = 'Normal'

$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +9.56%  '
$ws.Range('E50').Style = 'Normal'

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.220'
$ws.Range('D51').Style = 'Normal'

$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +0.35%  '
$ws.Range('E51').Style = 'Normal'
